$d = $word.ActiveDocument

# Find the paragraph holding "Ver no Jupiter Salvar em pdf Salvar em docx"
# and the paragraph holding the "© 2020 ..." footer line, then remove both
# of them together with the blank paragraph that immediately precedes the
# first one (mirroring the site-wide footer/navigation block that was
# dropped from the generated page).

$jupiterIndex = -1
$copyrightIndex = -1

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text
    if ($jupiterIndex -lt 0 -and $text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $jupiterIndex = $i
    }
    if ($copyrightIndex -lt 0 -and $text -like "*Contact: luizeleno@usp.br*") {
        $copyrightIndex = $i
    }
}

if ($jupiterIndex -gt 0 -and $copyrightIndex -ge $jupiterIndex) {
    $startPara = $d.Paragraphs.Item($jupiterIndex - 1)
    $endPara = $d.Paragraphs.Item($copyrightIndex)
    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
